$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: becomes "A 65163-2020" record (was row 3), with tweaks ---
$ws.Range("A2").Value2 = "A 65163-2020"
$ws.Range("B2").Value2 = 44172
$ws.Range("C2").Value2 = 45202
$ws.Range("D2").Value2 = "KALMAR LÄN"
$ws.Range("E2").Value2 = "MÖRBYLÅNGA"
$ws.Range("G2").Value2 = 16.8
$ws.Range("H2").Value2 = 2
$ws.Range("I2").Value2 = 3
$ws.Range("J2").Value2 = 6
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = 0
$ws.Range("N2").Value2 = 0
$ws.Range("O2").Value2 = 7
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 10
$ws.Range("R2").Value2 = "Liten diskröksvamp`r`nDvärgjordstjärna`r`nFyrflikig jordstjärna`r`nFågelarv`r`nKlibbveronika`r`nMindre hackspett`r`nRakhorndyvel`r`nEkoxe`r`nMurgröna`r`nRödbrun jordstjärna"
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/artfynd/A 65163-2020.xlsx", "A 65163-2020")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/kartor/A 65163-2020.png", "A 65163-2020")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/klagomål/A 65163-2020.docx", "A 65163-2020")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/klagomålsmail/A 65163-2020.docx", "A 65163-2020")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/tillsyn/A 65163-2020.docx", "A 65163-2020")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/tillsynsmail/A 65163-2020.docx", "A 65163-2020")'

# --- Row 3: becomes "A 66041-2018" record (was row 2); only date C changes ---
$ws.Range("A3").Value2 = "A 66041-2018"
$ws.Range("B3").Value2 = 43434
$ws.Range("C3").Value2 = 45202
$ws.Range("D3").Value2 = "KALMAR LÄN"
$ws.Range("E3").Value2 = "MÖRBYLÅNGA"
$ws.Range("G3").Value2 = 3.4
$ws.Range("H3").Value2 = 2
$ws.Range("I3").Value2 = 2
$ws.Range("J3").Value2 = 3
$ws.Range("K3").Value2 = 0
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 2
$ws.Range("N3").Value2 = 0
$ws.Range("O3").Value2 = 6
$ws.Range("P3").Value2 = 3
$ws.Range("Q3").Value2 = 9
$ws.Range("R3").Value2 = "Lundalm`r`nSkogsalm`r`nAsk`r`nDesmeknopp`r`nHårig jordstjärna`r`nÄngsskära`r`nTvåblad`r`nUnderviol`r`nSankt pers nycklar"
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/artfynd/A 66041-2018.xlsx", "A 66041-2018")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/kartor/A 66041-2018.png", "A 66041-2018")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/klagomål/A 66041-2018.docx", "A 66041-2018")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/klagomålsmail/A 66041-2018.docx", "A 66041-2018")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/tillsyn/A 66041-2018.docx", "A 66041-2018")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MORBYLANGA/tillsynsmail/A 66041-2018.docx", "A 66041-2018")'

# Row heights are unchanged by the source edit (still 15pt); the content
# rewrite above would otherwise trigger an autofit on the wrapped R column,
# so pin them back explicitly.
$ws.Rows("2").RowHeight = 15
$ws.Rows("3").RowHeight = 15

# --- Rows 4-52: only "Förändrad" date (column C) changes from 45192 to 45202 ---
for ($r = 4; $r -le 52; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45202
}

